# Daily attendance processing - swap the order of "System" and the
# recorder's email address in the "Recorded By" column (G) so that
# "System" is listed first, e.g. "dnasr281@gmail.com, System" becomes
# "System, dnasr281@gmail.com".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $value = $cell.Text

    if ($value -eq "dnasr281@gmail.com, System") {
        $cell.Value = "System, dnasr281@gmail.com"
    }
    elseif ($value -eq "admin@admin.com, System") {
        $cell.Value = "System, admin@admin.com"
    }
}
